$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 85 (shifts old rows 85-177 down to 86-178)
$ws.Rows(85).Insert()

# Populate the newly inserted row 85 with the new data record
$ws.Cells.Item(85, 1).Value = 8
$ws.Cells.Item(85, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(85, 3).Value = "Coquimbo"
$ws.Cells.Item(85, 4).Value = 44880
$ws.Cells.Item(85, 5).Value = 4
$ws.Cells.Item(85, 6).Value = 100112001
$ws.Cells.Item(85, 7).Value = "Berenjena"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 440
$ws.Cells.Item(85, 11).Value = 10000
$ws.Cells.Item(85, 12).Value = 11000
$ws.Cells.Item(85, 13).Value = 10500
$ws.Cells.Item(85, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(85, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(85, 16).Value = 262
$ws.Cells.Item(85, 17).Value = 40
$ws.Cells.Item(85, 18).Value = "Hortaliza"
